$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.990.69"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "2.613.10"
$ws.Range("E3").Value = "  -1.58%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'587.08"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").Value = "'165.19"
$ws.Range("E6").Value = "  -2.00%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.530"
$ws.Range("E8").Value = "  -2.57%  "

$ws.Range("D9").Value = "2.612.77"
$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("E10").Value = "  -4.22%  "

$ws.Range("E11").Value = "  +0.64%  "

$ws.Range("D12").Value = "'0.367"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "'5.20"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").Value = "'27.21"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("D15").Value = "3.094.40"
$ws.Range("E15").Value = "  -1.39%  "

$ws.Range("E16").Value = "  -2.86%  "

$ws.Range("D17").Value = "66.966.66"
$ws.Range("E17").Value = "  -0.87%  "

$ws.Range("D18").Value = "2.618.97"
$ws.Range("E18").Value = "  -1.53%  "

$ws.Range("D19").Value = "'11.69"
$ws.Range("E19").Value = "  -3.15%  "

$ws.Range("D20").Value = "'7.78"
$ws.Range("E20").Value = "  -6.52%  "

$ws.Range("D21").Value = "'355.15"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  -3.25%  "

$ws.Range("D23").Value = "'4.64"
$ws.Range("E23").Value = "  -3.50%  "

$ws.Range("D24").Value = "'10.54"
$ws.Range("E24").Value = "  -4.45%  "

$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").Value = "'1.91"
$ws.Range("E26").Value = "  -4.98%  "

$ws.Range("D27").Value = "'69.40"
$ws.Range("E27").Value = "  -2.30%  "

$ws.Range("D28").Value = "2.757.94"
$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "0.0₃0993"
$ws.Range("E30").Value = "  -3.55%  "

$ws.Range("D31").Value = "'542.00"
$ws.Range("E31").Value = "  -2.97%  "

$ws.Range("D32").Value = "'8.11"
$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("E33").Value = "  -4.46%  "

$ws.Range("E34").Value = "  -3.64%  "

$ws.Range("D35").Value = "'0.133"
$ws.Range("E35").Value = "  -0.71%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  -4.72%  "

$ws.Range("D38").Value = "'158.92"
$ws.Range("E38").Value = "  +0.34%  "

$ws.Range("D39").Value = "'18.91"
$ws.Range("E39").Value = "  -2.73%  "

$ws.Range("D40").Value = "'0.363"
$ws.Range("E40").Value = "  -2.68%  "

$ws.Range("D41").Value = "'18.24"
$ws.Range("E41").Value = "  +1.71%  "

$ws.Range("D42").Value = "'1.80"
$ws.Range("E42").Value = "  -1.86%  "

$ws.Range("D43").Value = "'5.12"
$ws.Range("E43").Value = "  -3.60%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").Value = "'2.41"
$ws.Range("E45").Value = "  -5.53%  "

$ws.Range("D46").Value = "0.0₆0297"
$ws.Range("E46").Value = "  -1.17%  "

$ws.Range("D47").Value = "'0.577"
$ws.Range("E47").Value = "  -3.76%  "

$ws.Range("D48").Value = "'150.64"
$ws.Range("E48").Value = "  -2.69%  "

$ws.Range("D49").Value = "'3.76"
$ws.Range("E49").Value = "  -3.45%  "

$ws.Range("E50").Value = "  -2.01%  "

$ws.Range("E51").Value = "  -1.80%  "
